# edit.ps1
# Applies the cryptos-list data refresh commit
# ("Updated cryptos list on Mon Aug 12 02:58:34 UTC 2024 with GitHub Actions").
#
# Only the "Price" (column D) and "Volume(1h)" (column E) cells change; every
# other cell (A/B/C columns, headers, formatting) is left untouched.
#
# Column D values are free-form text (e.g. "58.485.98", "6.18", "0.0₃0787") using
# "." as a thousands/decimal glyph rather than real numbers, so we prefix them with
# a leading apostrophe to force Excel to store them as text instead of silently
# re-parsing number-looking strings (like "6.18") into floating point values. The
# Style reset afterwards clears the quote-prefix flag so no stray formatting is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''58.485.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.04%  '
$ws.Range("D3").Value = '''2.537.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.11%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''507.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.20%  '
$ws.Range("D6").Value = '''144.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.01%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.562'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.01%  '
$ws.Range("D9").Value = '''2.538.28'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("D10").Value = '''6.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.41%  '
$ws.Range("E11").Value = '  -6.77%  '
$ws.Range("E12").Value = '  -4.79%  '
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("D14").Value = '''2.979.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.20%  '
$ws.Range("D15").Value = '''58.436.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.09%  '
$ws.Range("D16").Value = '''20.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.75%  '
$ws.Range("E17").Value = '  -6.19%  '
$ws.Range("D18").Value = '''2.536.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = '  -5.32%  '
$ws.Range("D20").Value = '''334.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.53%  '
$ws.Range("D21").Value = '''10.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.12%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("E23").Value = '  -4.67%  '
$ws.Range("D24").Value = '''60.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.25%  '
$ws.Range("E25").Value = '  -5.12%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("E27").Value = '  -5.27%  '
$ws.Range("D28").Value = '''2.646.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.27%  '
$ws.Range("D29").Value = '''0.0₃0787'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.08%  '
$ws.Range("D30").Value = '''6.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.99%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").Value = '''149.86'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("D33").Value = '''5.84'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.50%  '
$ws.Range("D34").Value = '''18.50'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.92%  '
$ws.Range("E35").Value = '  -5.25%  '
$ws.Range("D36").Value = '''0.942'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.57%  '
$ws.Range("D37").Value = '''3.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.07%  '
$ws.Range("E38").Value = '  -7.64%  '
$ws.Range("D39").Value = '''36.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("D40").Value = '''0.825'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -11.45%  '
$ws.Range("E41").Value = '  -6.49%  '
$ws.Range("D42").Value = '''284.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.79%  '
$ws.Range("E43").Value = '  -7.34%  '
$ws.Range("D44").Value = '''0.0996'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.23%  '
$ws.Range("D45").Value = '''0.996'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("E46").Value = '  -5.77%  '
$ws.Range("E47").Value = '  -4.96%  '
$ws.Range("D48").Value = '''18.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.41%  '
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("E50").Value = '  -5.18%  '
$ws.Range("D51").Value = '''4.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.96%  '
